$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.141.78"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.64"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.92"
$ws.Range("E5").Value = "  -2.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5219"
$ws.Range("E6").Value = "  -1.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2626"
$ws.Range("E8").Value = "  -2.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06330"
$ws.Range("E9").Value = "  -0.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.18"
$ws.Range("E10").Value = "  -1.36%  "

$ws.Range("E11").Value = "  -1.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.665.61"
$ws.Range("E12").Value = "  -1.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.431"
$ws.Range("E13").Value = "  -2.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5496"
$ws.Range("E14").Value = "  -4.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008024"
$ws.Range("E15").Value = "  -2.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.45"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.165.66"
$ws.Range("E17").Value = "  -0.30%  "

$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.748"
$ws.Range("E19").Value = "  -2.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.32"
$ws.Range("E20").Value = "  -1.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.31"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.229"
$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.82"
$ws.Range("E24").Value = "  +0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1239"
$ws.Range("E25").Value = "  -1.82%  "

$ws.Range("E26").Value = "  -3.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.82"
$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06315"
$ws.Range("E28").Value = "  -1.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.355"
$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  -2.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.521"
$ws.Range("E31").Value = "  -1.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.417"
$ws.Range("E32").Value = "  -4.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.649"
$ws.Range("E33").Value = "  -1.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.006"
$ws.Range("E34").Value = "  -1.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6025"
$ws.Range("E35").Value = "  -1.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.399"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.763"
$ws.Range("E37").Value = "  +0.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.117.75"
$ws.Range("E38").Value = "  +1.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.086"
$ws.Range("E39").Value = "  -1.52%  "

$ws.Range("E40").Value = "  -0.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8645"
$ws.Range("E41").Value = "  -1.99%  "

$ws.Range("E42").Value = "  -0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.50"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.822.07"
$ws.Range("E44").Value = "  -0.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("E45").Value = "  -0.52%  "

$ws.Range("E46").Value = "  -3.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.061"
$ws.Range("E48").Value = "  -0.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05237"
$ws.Range("E49").Value = "  -0.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4241"
$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.922"
$ws.Range("E51").Value = "  -1.56%  "
